$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.014.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.638.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.75"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.665.36"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.24"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.102.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.008.93"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.92"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.660.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.23"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.36"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.987"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.23%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.01"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.08"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.970"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.65"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.841"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.71"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.40"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.93"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.612"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0988"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.994"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0530"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.30"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0230"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.985.63"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.67"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.71%  "
